$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values for column G, rows 2-28.
# Old values were previously computed from the "Strike#" source column;
# these are regenerated values based on the actual K (strikeouts) stat.
$kValues = @{
    2  = 7
    3  = 6
    4  = 6
    5  = 3
    6  = 6
    7  = 6
    8  = 5
    9  = 7
    10 = 5
    11 = 4
    12 = 0
    13 = 4
    14 = 9
    15 = 8
    16 = 5
    17 = 4
    18 = 11
    19 = 4
    20 = 3
    21 = 5
    22 = 3
    23 = 3
    24 = 10
    25 = 2
    26 = 1
    27 = 5
    28 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
